# Updated cryptos list on Sun Sep 22 23:26:09 UTC 2024 with GitHub Actions
#
# Applies per-row Price (column D) and Volume(1h) (column E) refresh values,
# plus two row-content swaps (rows 27<->28 and rows 32<->33) where the
# ranking reordered two adjacent coins - matching the target OOXML diff.
#
# Every column-D write is preceded by forcing the cell's NumberFormat to
# text ("@") because these Price values ("63.558.58", "0.587", "1.00", ...)
# are stored as plain text in the workbook (not numbers), and Excel would
# otherwise auto-coerce numeric-looking strings into real numbers
# (dropping significant trailing zeros, e.g. "1.00" -> 1, "0.0960" -> 0.096).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.558.58'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.577.94'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '586.06'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.75'
$ws.Range('E6').Value = '  -2.79%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  -1.84%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.106'
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.351'
$ws.Range('E12').Value = '  -1.81%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.33'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.041.60'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '63.480.04'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('E16').Value = '  -2.44%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.579.95'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.13'
$ws.Range('E18').Value = '  -2.59%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '342.33'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.31'
$ws.Range('E20').Value = '  -3.30%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.64'
$ws.Range('E21').Value = '  -3.79%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.26'
$ws.Range('E23').Value = '  +2.20%  '
$ws.Range('E24').Value = '  +6.15%  '
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('E26').Value = '  -3.96%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.95'
$ws.Range('E28').Value = '  -3.73%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.24'
$ws.Range('E29').Value = '  -3.69%  '
$ws.Range('E30').Value = '  -2.64%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '471.66'
$ws.Range('E31').Value = '  +0.96%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0₃0801'
$ws.Range('E32').Value = '  -3.68%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.68'
$ws.Range('E33').Value = '  +2.24%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '176.56'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.89'
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.72'
$ws.Range('E40').Value = '  -2.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '160.27'
$ws.Range('E41').Value = '  +5.33%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '39.98'
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.71'
$ws.Range('E43').Value = '  -3.82%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.74'
$ws.Range('E44').Value = '  +2.72%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.631'
$ws.Range('E45').Value = '  +2.21%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0538'
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0960'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0237'
$ws.Range('E48').Value = '  -2.01%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '18.18'
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '11.37'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.69'
$ws.Range('E51').Value = '  -4.47%  '
